$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "33.33 / 100.0"
$ws.Range("C2").Value = "student_display_name_with_only_last_name: failed`nstudent_display_full_name: failed"
$ws.Range("D2").Value = "33.33/100.0"

$ws.Range("B3").Value = "100.0 / 100.0"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "100.0/100.0"

$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
